# Weekly update: insert a new price record for "Acelga" (Vega Modelo de
# Temuco) at the top of the data block (row 181), pushing the existing
# historical rows (181-218) down by one (to 182-219).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 181; rows 181..218 shift down to 182..219 and the
# sheet's used range/dimension grows to A1:R219 automatically.
$ws.Rows.Item(181).Insert()

# Populate the newly-inserted row 181 with the new weekly observation.
$ws.Cells.Item(181, 1).Value  = 10
$ws.Cells.Item(181, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(181, 3).Value  = "La Araucanía"
$ws.Cells.Item(181, 4).Value  = 44511
$ws.Cells.Item(181, 5).Value  = 9
$ws.Cells.Item(181, 6).Value  = 100112009
$ws.Cells.Item(181, 7).Value  = "Acelga"
$ws.Cells.Item(181, 8).Value  = "Sin especificar"
$ws.Cells.Item(181, 9).Value  = "Primera"
$ws.Cells.Item(181, 10).Value = 60
$ws.Cells.Item(181, 11).Value = 8000
$ws.Cells.Item(181, 12).Value = 9000
$ws.Cells.Item(181, 13).Value = 8500
$ws.Cells.Item(181, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(181, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(181, 16).Value = 708
$ws.Cells.Item(181, 17).Value = 12
$ws.Cells.Item(181, 18).Value = "Hortaliza"
